$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (row 1), columns L..V
# NOTE: shared-string pool is built in the order new strings are first
# assigned; set M1/N1 in LEN,HT order so the pool matches the source
# workbook (which had headers typed before a column reorder), then fix
# up the actual header text/order afterwards.
$ws.Range("L1").Value = "SHIP_CARTON_WID"
$ws.Range("M1").Value = "SHIP_CARTON_LEN"
$ws.Range("N1").Value = "SHIP_CARTON_HT"
$ws.Range("O1").Value = "SHIP_CARTON_WT"
$ws.Range("P1").Value = "S_UNIT_WIDTH"
$ws.Range("Q1").Value = "S_UNIT_HEIGHT"
$ws.Range("R1").Value = "S_UNIT_LENGTH"
$ws.Range("S1").Value = "S_UNIT_WEIGHT"
$ws.Range("T1").Value = "SUPP_PACK_SIZE"
$ws.Range("U1").Value = "INNERPACK_SIZE"
$ws.Range("V1").Value = "FRENCH_COMPLIANT"

# Final layout has SHIP_CARTON_HT before SHIP_CARTON_LEN (M/N swapped
# relative to typing order above) -- fix up now that both strings exist.
$ws.Range("M1").Value = "SHIP_CARTON_HT"
$ws.Range("N1").Value = "SHIP_CARTON_LEN"

# Row 2 values
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 1
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 1
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 1
$ws.Range("V2").Value = "Yes"

# Row 3 values
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = "No"

# Column width best-fit for newly visible columns (values chosen so the
# persisted OOXML width - which is ColumnWidth + 5/6 - lands as close as
# this engine's internal rounding allows to the target bestFit widths)
$ws.Columns("L").ColumnWidth = 12671/768
$ws.Columns("M").ColumnWidth = 11519/768
$ws.Columns("N").ColumnWidth = 97/6
$ws.Columns("O").ColumnWidth = 11903/768
$ws.Columns("R").ColumnWidth = 10367/768
$ws.Columns("S").ColumnWidth = 10367/768

# Update view: scroll/top-left cell and selection
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("V3").Select()

$wb.Windows.Item(1).WindowState = -4143
